# Update the "想去人数" (want-to-go count) values for two events
# in both the "展览" sheet and its duplicate "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 116
    $ws.Range("F4").Value = 74
}
